# Common: Added more cool stuff
# Adds new "base" (báze) translation keys to the Import sheet, mirroring
# the existing "liquid" translation block pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New key/value pairs (column B = translation key, column C = Czech value)
$rows = @(
    @("lab.base.tooltip.create", "Nová báze"),
    @("lab.base.tooltip.create", "Nová báze"),
    @("lab.base.create.title", "Nová báze"),
    @("lab.base.create.subtitle", "Báze je užitečná v mixech."),
    @("lab.base.name.label", "Název báze"),
    @("lab.base.vendorId.label", "Výrobce"),
    @("lab.base.pg.label", "PG"),
    @("lab.base.vg.label", "VG"),
    @("lab.base.create.submit", "Vytvořit bázi"),
    @("lab.base.create.success", "Báze [{{data.name}}] byla uložena.")
)

$lastRow = 344
$startRow = $lastRow + 1
$endRow = $startRow + $rows.Length - 1

# Copy the formatting (style) of the last existing data row down across the
# newly added rows so they keep the same "import" cell style.
$ws.Range("A344:C344").Copy()
$ws.Range("A" + $startRow + ":C" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value2 = "cs"
    $ws.Cells.Item($r, 2).Value2 = $row[0]
    $ws.Cells.Item($r, 3).Value2 = $row[1]
    $r = $r + 1
}

# Update the view so it mirrors the scrolled/selected state of the edited sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 337
$ws.Range("B353").Select()
